# Edit script: insert new job row at top, shift existing rows down,
# update timestamps, refresh hyperlinks, and widen column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting the existing data rows down by one.
$ws.Rows.Item(2).Insert()

# Row data (rows 2-10 of the final sheet), columns A-H.
$rowData = @(
    @("2025-12-12 12:39:14", "初回 スポーツクラブ コスパ自動予約bot開発(playwight/Python)", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5452614", 368, "🔥Python ★bot ◆開発"),
    @("2025-12-12 12:39:14", "AIオートメーションエンジニア", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5452520", 303, "🔥AI,Ai"),
    @("2025-12-12 12:39:14", "【Flutterエンジニア募集】Androidアプリ開発のパートナーを探しています", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5452211", 100, "◆開発 ◇アプリ"),
    @("2025-12-12 12:39:14", "【再掲】基幹システム入替に伴うBIツール環境の再構築(Microsoft Power BI)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5452367", 88, "◆ツール"),
    @("2025-12-12 12:39:14", "Amazonの購入アカウントから必要な情報のスクレイピング→スプレッドシートに記入をしたい。", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5452210", 40, "◆スクレイピング"),
    @("2025-12-12 12:39:14", "【R/Shiny】高齢者評価アプリ 機能追加・UI改修依頼", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5452159", 38, "◇アプリ"),
    @("2025-12-12 12:39:14", "【小規模・短納期・急募】アプリMatrixifyを用いたデータ移行検証・マッピング担当募集", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5451926", 33, "◇アプリ"),
    @("2025-12-12 12:39:14", "注目 限定公開 PR 限定公開の仕事", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5450323", 13, $null),
    @("2025-12-12 12:39:14", "Xの運用代行", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5451931", 10, $null),
)

# Write all cell values for rows 2-10.
$r = 2
foreach ($row in $rowData) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $val = $row[$i]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $i + 1).Value = $val
        }
    }
    $r++
}

# Rebuild the hyperlinks for column F (URLs), in row order, so relationship
# ids rId1..rId9 line up with rows 2..10 top to bottom.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 10; $r++) {
    $url = $ws.Cells.Item($r, 6).Value()
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $url)
}

# Widen column H (skill summary) from 12 to 18 characters.
# (engine adds a constant ~0.8333 padding to ColumnWidth assignments, so
#  compensate to land on an exact width of 18 in the saved file)
$ws.Columns.Item(8).ColumnWidth = (18 - 5/6)

